$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report's single data row (row 4) is cleared out (no transactions in
# this period) and its counter reset to 0. Row 5, which held the second
# data row, is removed entirely - this shifts the old totals row (6) up to
# become row 5, and the old footer row (7) up to become row 6.
$ws.Range("B4:G4").ClearContents()
$ws.Range("H4:K4").ClearContents()
$ws.Range("L4:M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("A4").Value = 0

$ws.Rows.Item(5).Delete()

# The totals row (now row 5) loses its computed total along with the data.
$ws.Range("K5:N5").ClearContents()

Write-Output "done"
